$d = $word.ActiveDocument

# --- Edit 1: first paragraph -------------------------------------------
# "This is a Microsoft word document." gains two trailing spaces, then a
# red-colored parenthetical is appended in three separate runs, matching
# how Word splits runs for text typed/spell-checked incrementally.
$p1 = $d.Paragraphs.Item(1)

# Append two trailing spaces to the existing (black) run.
$p1.Range.InsertAfter("  ")

# Red run 1: "(This is a change – Ve"
$r1 = $d.Range($p1.Range.End - 1, $p1.Range.End - 1)
$r1.InsertAfter("(This is a change " + [char]0x2013 + " Ve")
$r1.Font.Color = 255

# Red run 2: "rsion for main branch"
$r2 = $d.Range($p1.Range.End - 1, $p1.Range.End - 1)
$r2.InsertAfter("rsion for main branch")
$r2.Font.Color = 255

# Red run 3: ")"
$r3 = $d.Range($p1.Range.End - 1, $p1.Range.End - 1)
$r3.InsertAfter(")")
$r3.Font.Color = 255

# --- Edit 2: remove the trailing "ank God almighty..." paragraph -------
# This paragraph (the very last one, using the NormalWeb style) is
# deleted outright, leaving "Shall be lifted—nevermore!" as the final
# paragraph of the document.
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastIndex)
$lastPara.Range.Delete()
